$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 74 -> Problem ID 73: "Count Number of Nice Sub Arrays"
$ws.Range("B74").Value = "Prefix Sum"
$ws.Range("C74").Value = "Count Number of Nice Sub Arrays"
$ws.Range("D74").Value = "Medium"
$ws.Range("E74").Value = "Done"
$ws.Range("F73").Copy()
$ws.Range("F74").PasteSpecial(-4122)
$ws.Range("F74").Value = 45914
$ws.Range("G74").Value = "O(n)"
$ws.Range("H74").Value = "O(n)"
$ws.Range("I74").Value = "Prefix Sum"

# Row 75 -> Problem ID 74: "Minimum Operations to Reduce X to Zero"
$ws.Range("B75").Value = "Sliding Window"
$ws.Range("C75").Value = "Minimum Operations to Reduce X to Zero"
$ws.Range("D75").Value = "Medium"
$ws.Range("E75").Value = "Done"
$ws.Range("F73").Copy()
$ws.Range("F75").PasteSpecial(-4122)
$ws.Range("F75").Value = 45914
$ws.Range("G75").Value = "O(n)"
$ws.Range("H75").Value = "O(1)"
$ws.Range("I75").Value = "Sliding Window"

# Match the new selection recorded in the workbook after the edit
$ws.Range("C76").Select()
